# Login with Invalid Username Test
#
# The DataSet sheet's row 9 ("LoginWithInvalidUser") used to hold a
# hyperlinked sample e-mail address (abv@abv.bg) in B9. The test value is
# replaced with a plain (non e-mail) invalid username "sgdfhh", so the
# hyperlink + hyperlink styling on B9 needs to go away and the cell reverts
# to the ordinary (non-hyperlink) look used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the hyperlink that lives on B9 ---------------------------------
# The object model here only offers a sheet-wide Hyperlinks.Delete(), so
# record every other hyperlinked cell + target first, wipe all hyperlinks,
# then recreate every one of them except the one that used to sit on B9.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:abv@abv.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 2), "mailto:Ilian@mail.bg", "", "", "Ilian@mail.bg") | Out-Null
$ws.Cells.Item(6, 2).Value = "rosen"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 2), "mailto:abv@abv.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 2), "mailto:abv@abv.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 2), "mailto:abv@abv.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 2), "mailto:abv@abv.bg") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 2), "mailto:abv@abv.bg") | Out-Null

# --- B9: new value, plain (non-hyperlink) formatting ------------------------
$b9 = $ws.Cells.Item(9, 2)
$b9.Clear()
$b9.Value = "sgdfhh"
$b9.NumberFormat = "@"
$b9.HorizontalAlignment = -4108

# --- Row 9 no longer needs a custom row height ------------------------------
$ws.Rows.Item(9).AutoFit()

# --- Selection moves to B9 ---------------------------------------------------
$ws.Range("B9").Select()
